# Market-data refresh: refreshes cached currentAveragePrice / LevePrice /
# LeveProfit figures (columns H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# crafting-leve sheets, as produced by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1888.8889
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 2375
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 2375
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -2725
$ws.Range("H76").Value = 4385.7144
$ws.Range("J76").Value = 4316.6665
$ws.Range("L76").Value = 4316.6665
$ws.Range("N76").Value = -4946.6665
$ws.Range("H79").Value = 4385.7144
$ws.Range("J79").Value = 4316.6665
$ws.Range("L79").Value = 4316.6665
$ws.Range("N79").Value = -6500.6665
$ws.Range("H112").Value = 2003.7188
$ws.Range("J112").Value = 2117.862
$ws.Range("L112").Value = 6353.586
$ws.Range("N112").Value = -8569.585999999999
$ws.Range("H132").Value = 1324.2603
$ws.Range("I132").Value = 1137.0454
$ws.Range("J132").Value = 3089.4285
$ws.Range("K132").Value = 3411.1362
$ws.Range("L132").Value = 9268.2855
$ws.Range("M132").Value = -881.1361999999999
$ws.Range("N132").Value = -14328.2855
$ws.Range("H133").Value = 49816
$ws.Range("J133").Value = 49816
$ws.Range("L133").Value = 49816
$ws.Range("N133").Value = -59936
$ws.Range("H138").Value = 2956.82
$ws.Range("I138").Value = 2020.2
$ws.Range("J138").Value = 3581.2334
$ws.Range("K138").Value = 6060.6
$ws.Range("L138").Value = 10743.7002
$ws.Range("M138").Value = -920.6000000000004
$ws.Range("N138").Value = -21023.7002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1762.6744
$ws.Range("I61").Value = 1682.9474
$ws.Range("J61").Value = 2368.6
$ws.Range("K61").Value = 1682.9474
$ws.Range("L61").Value = 2368.6
$ws.Range("M61").Value = -1470.9474
$ws.Range("N61").Value = -2792.6
$ws.Range("H74").Value = 1284.7358
$ws.Range("I74").Value = 1123.579
$ws.Range("J74").Value = 1693
$ws.Range("K74").Value = 1123.579
$ws.Range("L74").Value = 1693
$ws.Range("M74").Value = -249.579
$ws.Range("N74").Value = -3441
$ws.Range("H77").Value = 1284.7358
$ws.Range("I77").Value = 1123.579
$ws.Range("J77").Value = 1693
$ws.Range("K77").Value = 5617.895
$ws.Range("L77").Value = 8465
$ws.Range("M77").Value = -1249.895
$ws.Range("N77").Value = -17201
$ws.Range("H101").Value = 52734.668
$ws.Range("J101").Value = 52734.668
$ws.Range("L101").Value = 52734.668
$ws.Range("N101").Value = -59224.668
$ws.Range("H122").Value = 3056.6
$ws.Range("J122").Value = 2304.6667
$ws.Range("L122").Value = 6914.000100000001
$ws.Range("N122").Value = -11814.0001
$ws.Range("H132").Value = 393520.9
$ws.Range("I132").Value = 477169.66
$ws.Range("K132").Value = 1431508.98
$ws.Range("M132").Value = -1428978.98
$ws.Range("H136").Value = 1762.6744
$ws.Range("I136").Value = 1682.9474
$ws.Range("J136").Value = 2368.6
$ws.Range("K136").Value = 5048.8422
$ws.Range("L136").Value = 7105.799999999999
$ws.Range("M136").Value = -2498.8422
$ws.Range("N136").Value = -12205.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 590.25
$ws.Range("I64").Value = 685.3333
$ws.Range("J64").Value = 305
$ws.Range("K64").Value = 685.3333
$ws.Range("L64").Value = 305
$ws.Range("M64").Value = -460.3333
$ws.Range("N64").Value = -755
$ws.Range("H67").Value = 590.25
$ws.Range("I67").Value = 685.3333
$ws.Range("J67").Value = 305
$ws.Range("K67").Value = 685.3333
$ws.Range("L67").Value = 305
$ws.Range("M67").Value = 94.66669999999999
$ws.Range("N67").Value = -1865
$ws.Range("H105").Value = 3162.4285
$ws.Range("I105").Value = 2811.7222
$ws.Range("K105").Value = 2811.7222
$ws.Range("M105").Value = -1064.7222
$ws.Range("H134").Value = 291816.66
$ws.Range("I134").Value = 340497.56
$ws.Range("J134").Value = 4599.4
$ws.Range("K134").Value = 1021492.68
$ws.Range("L134").Value = 13798.2
$ws.Range("M134").Value = -1018957.68
$ws.Range("N134").Value = -18868.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3038.525
$ws.Range("I31").Value = 2217.1
$ws.Range("K31").Value = 2217.1
$ws.Range("M31").Value = -1922.1
$ws.Range("H34").Value = 3038.525
$ws.Range("I34").Value = 2217.1
$ws.Range("K34").Value = 2217.1
$ws.Range("M34").Value = -2015.1
$ws.Range("H51").Value = 32500
$ws.Range("J51").Value = 32500
$ws.Range("L51").Value = 32500
$ws.Range("N51").Value = -33972
$ws.Range("H58").Value = 1324066.1
$ws.Range("I58").Value = 2180096
$ws.Range("J58").Value = 1110.7273
$ws.Range("K58").Value = 2180096
$ws.Range("L58").Value = 1110.7273
$ws.Range("M58").Value = -2179893
$ws.Range("N58").Value = -1516.7273
$ws.Range("H61").Value = 32500
$ws.Range("J61").Value = 32500
$ws.Range("L61").Value = 32500
$ws.Range("N61").Value = -33196
$ws.Range("H132").Value = 348431.25
$ws.Range("I132").Value = 616132.8
$ws.Range("J132").Value = 1993.9412
$ws.Range("K132").Value = 1848398.4
$ws.Range("L132").Value = 5981.8236
$ws.Range("M132").Value = -1845868.4
$ws.Range("N132").Value = -11041.8236
$ws.Range("H134").Value = 2192.6365
$ws.Range("I134").Value = 2249.875
$ws.Range("J134").Value = 2040
$ws.Range("K134").Value = 6749.625
$ws.Range("L134").Value = 6120
$ws.Range("M134").Value = -4214.625
$ws.Range("N134").Value = -11190
$ws.Range("H136").Value = 1324066.1
$ws.Range("I136").Value = 2180096
$ws.Range("J136").Value = 1110.7273
$ws.Range("K136").Value = 6540288
$ws.Range("L136").Value = 3332.1819
$ws.Range("M136").Value = -6537738
$ws.Range("N136").Value = -8432.1819
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19611516
$ws.Range("J131").Value = 27028994
$ws.Range("L131").Value = 81086982
$ws.Range("N131").Value = -81097062
$ws.Range("H138").Value = 1259.4445
$ws.Range("I138").Value = 980.625
$ws.Range("J138").Value = 3490
$ws.Range("K138").Value = 2941.875
$ws.Range("L138").Value = 10470
$ws.Range("M138").Value = 2198.125
$ws.Range("N138").Value = -20750
$ws.Range("H139").Value = 1608.4
$ws.Range("I139").Value = 1168.9474
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 3506.8422
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 1633.1578
$ws.Range("N139").Value = -19280
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 4931.864
$ws.Range("I122").Value = 4815.6924
$ws.Range("J122").Value = 5099.6665
$ws.Range("K122").Value = 14447.0772
$ws.Range("L122").Value = 15298.9995
$ws.Range("M122").Value = -11997.0772
$ws.Range("N122").Value = -20198.9995
$ws.Range("H132").Value = 2312.2222
$ws.Range("I132").Value = 1976.25
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5928.75
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3398.75
$ws.Range("N132").Value = -20060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7692.4116
$ws.Range("I40").Value = 7897.1
$ws.Range("J40").Value = 7400
$ws.Range("K40").Value = 7897.1
$ws.Range("L40").Value = 7400
$ws.Range("M40").Value = -7761.1
$ws.Range("N40").Value = -7672
$ws.Range("H48").Value = 50000
$ws.Range("I48").Value = 50000
$ws.Range("K48").Value = 50000
$ws.Range("M48").Value = -49339
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H122").Value = 6444.7
$ws.Range("I122").Value = 6849.2144
$ws.Range("J122").Value = 5500.8335
$ws.Range("K122").Value = 20547.6432
$ws.Range("L122").Value = 16502.5005
$ws.Range("M122").Value = -18097.6432
$ws.Range("N122").Value = -21402.5005
$ws.Range("H132").Value = 6909.884
$ws.Range("I132").Value = 6771.206
$ws.Range("J132").Value = 7433.778
$ws.Range("K132").Value = 20313.618
$ws.Range("L132").Value = 22301.334
$ws.Range("M132").Value = -17783.618
$ws.Range("N132").Value = -27361.334
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 706.4545000000001
$ws.Range("I107").Value = 713.3333
$ws.Range("J107").Value = 698.2
$ws.Range("K107").Value = 2139.9999
$ws.Range("L107").Value = 2094.6
$ws.Range("M107").Value = -219.9998999999998
$ws.Range("N107").Value = -5934.6
$ws.Range("H122").Value = 1927.1428
$ws.Range("I122").Value = 1497.5
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4492.5
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2042.5
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 1606.7435
$ws.Range("I132").Value = 1345.0385
$ws.Range("J132").Value = 2130.1538
$ws.Range("K132").Value = 4035.1155
$ws.Range("L132").Value = 6390.4614
$ws.Range("M132").Value = -1505.1155
$ws.Range("N132").Value = -11450.4614
$ws.Range("H136").Value = 1889.0714
$ws.Range("I136").Value = 2156.423
$ws.Range("J136").Value = 1454.625
$ws.Range("K136").Value = 6469.268999999999
$ws.Range("L136").Value = 4363.875
$ws.Range("M136").Value = -3919.268999999999
$ws.Range("N136").Value = -9463.875
